$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "59.160.71"
$ws.Range("E2").Value = "  -2.74%  "

$ws.Range("D3").Value = "2.655.14"
$ws.Range("E3").Value = "  -1.37%  "

$ws.Range("E4").Value = "  -0.05%  "

$style_D5 = $ws.Range("D5").Style
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "523.81"
$ws.Range("D5").Style = $style_D5
$ws.Range("E5").Value = "  +0.18%  "

$style_D6 = $ws.Range("D6").Style
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "144.67"
$ws.Range("D6").Style = $style_D6
$ws.Range("E6").Value = "  -2.45%  "

$ws.Range("E7").Value = "  +0.31%  "

$ws.Range("E8").Value = "  -1.50%  "

$style_D9 = $ws.Range("D9").Style
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "7.05"
$ws.Range("D9").Style = $style_D9
$ws.Range("E9").Value = "  +9.49%  "

$ws.Range("E10").Value = "  -3.99%  "

$ws.Range("E11").Value = "  -2.55%  "

$style_D12 = $ws.Range("D12").Style
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.130"
$ws.Range("D12").Style = $style_D12
$ws.Range("E12").Value = "  +1.67%  "

$ws.Range("D13").Value = "3.119.70"
$ws.Range("E13").Value = "  -1.60%  "

$ws.Range("D14").Value = "59.172.14"
$ws.Range("E14").Value = "  -2.82%  "

$style_D15 = $ws.Range("D15").Style
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "21.17"
$ws.Range("D15").Style = $style_D15
$ws.Range("E15").Value = "  -1.66%  "

$ws.Range("E16").Value = "  -2.39%  "

$ws.Range("D17").Value = "2.659.67"
$ws.Range("E17").Value = "  -4.38%  "

$style_D18 = $ws.Range("D18").Style
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "340.33"
$ws.Range("D18").Style = $style_D18
$ws.Range("E18").Value = "  -4.31%  "

$style_D19 = $ws.Range("D19").Style
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "4.39"
$ws.Range("D19").Style = $style_D19
$ws.Range("E19").Value = "  -4.46%  "

$ws.Range("E20").Value = "  -2.06%  "

$ws.Range("E21").Value = "  +0.01%  "

$style_D22 = $ws.Range("D22").Style
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "0.997"
$ws.Range("D22").Style = $style_D22
$ws.Range("E22").Value = "  -0.36%  "

$style_D23 = $ws.Range("D23").Style
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "64.22"
$ws.Range("D23").Style = $style_D23
$ws.Range("E23").Value = "  +2.21%  "

$ws.Range("E24").Value = "  -2.77%  "

$ws.Range("E25").Value = "  -1.84%  "

$style_D26 = $ws.Range("D26").Style
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "0.998"
$ws.Range("D26").Style = $style_D26
$ws.Range("E26").Value = "  +0.50%  "

$ws.Range("E28").Value = "  -3.02%  "

$style_D29 = $ws.Range("D29").Style
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "6.69"
$ws.Range("D29").Style = $style_D29
$ws.Range("E29").Value = "  -2.08%  "

$style_D30 = $ws.Range("D30").Style
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "0.999"
$ws.Range("D30").Style = $style_D30
$ws.Range("E30").Value = "  +0.11%  "

$ws.Range("E31").Value = "  -0.39%  "

$ws.Range("E32").Value = "  -2.08%  "

$style_D33 = $ws.Range("D33").Style
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "149.29"
$ws.Range("D33").Style = $style_D33
$ws.Range("E33").Value = "  -0.60%  "

$style_D34 = $ws.Range("D34").Style
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "4.16"
$ws.Range("D34").Style = $style_D34
$ws.Range("E34").Value = "  -1.50%  "

$ws.Range("E35").Value = "  -2.91%  "

$style_D36 = $ws.Range("D36").Style
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.901"
$ws.Range("D36").Style = $style_D36
$ws.Range("E36").Value = "  -5.37%  "

$ws.Range("E37").Value = "  +0.38%  "

$ws.Range("B38").Value = "OKB"
$ws.Range("C38").Value = "https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb"
$style_D38 = $ws.Range("D38").Style
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "36.74"
$ws.Range("D38").Style = $style_D38
$ws.Range("E38").Value = "  +0.09%  "

$ws.Range("B39").Value = "Stacks"
$ws.Range("C39").Value = "https://coinranking.com/coin/mMPrMcB7+stacks-stx"
$style_D39 = $ws.Range("D39").Style
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "1.48"
$ws.Range("D39").Style = $style_D39
$ws.Range("E39").Value = "  -5.92%  "

$ws.Range("E40").Value = "  -4.17%  "

$ws.Range("E41").Value = "  +0.33%  "

$ws.Range("E42").Value = "  +0.42%  "

$style_D43 = $ws.Range("D43").Style
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "275.22"
$ws.Range("D43").Style = $style_D43
$ws.Range("E43").Value = "  -4.16%  "

$style_D44 = $ws.Range("D44").Style
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "19.86"
$ws.Range("D44").Style = $style_D44
$ws.Range("E44").Value = "  -1.14%  "

$ws.Range("E45").Value = "  -2.50%  "

$style_D46 = $ws.Range("D46").Style
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.0535"
$ws.Range("D46").Style = $style_D46
$ws.Range("E46").Value = "  -1.42%  "

$ws.Range("B48").Value = "Maker"
$ws.Range("C48").Value = "https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr"
$ws.Range("D48").Value = "2.031.85"
$ws.Range("E48").Value = "  -5.59%  "

$ws.Range("B49").Value = "RenderToken"
$ws.Range("C49").Value = "https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr"
$style_D49 = $ws.Range("D49").Style
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "4.77"
$ws.Range("D49").Style = $style_D49
$ws.Range("E49").Value = "  -2.97%  "

$ws.Range("E50").Value = "  -2.83%  "

$style_D51 = $ws.Range("D51").Style
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "18.95"
$ws.Range("D51").Style = $style_D51
$ws.Range("E51").Value = "  -1.55%  "
